# Fruta / hortaliza, semanal
# Insert the new week's record at the top of the data block (row 9),
# pushing the existing rows down by one. The oldest existing row
# (previously row 19) ends up duplicated as the new last row (row 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 9:19 down to 10:20, creating a blank row 9.
$ws.Rows("9:9").Insert()

# Populate the new row 9 with this week's data.
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C9").Value = "Los Lagos"
$ws.Range("D9").Value = 44771
$ws.Range("E9").Value = 10
$ws.Range("F9").Value = 100112013
$ws.Range("G9").Value = "Alcachofa"
$ws.Range("H9").Value = "Madrigal"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 90
$ws.Range("K9").Value = 16000
$ws.Range("L9").Value = 16000
$ws.Range("M9").Value = 16000
$ws.Range("N9").Value = "$/caja 40 unidades"
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 400
$ws.Range("Q9").Value = 40
$ws.Range("R9").Value = "Hortaliza"
